$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.96699333333333
$ws.Range("H2").Value = 38.90098
$ws.Range("I2").Value = 0.248939824633031
$ws.Range("J2").Value = 0.248939824633031
$ws.Range("M2").Value = 0.3667156666666667
$ws.Range("N2").Value = 1.100147
$ws.Range("O2").Value = 0.001893484040582144
$ws.Range("P2").Value = 0.001893484040582144
$ws.Range("Q2").Value = 4.755199604895555
$ws.Range("R2").Value = 42.79679644405999
$ws.Range("S2").Value = 0.0004713635850079618
$ws.Range("T2").Value = 0.0004713635850079618

$ws.Range("G3").Value = 12.96699333333333
$ws.Range("H3").Value = 38.90098
$ws.Range("I3").Value = 0.248939824633031
$ws.Range("J3").Value = 0.248939824633031
$ws.Range("N3").Value = 5.559651000000001
$ws.Range("O3").Value = 0.009568821657202681
$ws.Range("P3").Value = 0.009568821657202679
$ws.Range("Q3").Value = 24.03065248422
$ws.Range("R3").Value = 216.27587235798
$ws.Range("S3").Value = 0.002382060785288784
$ws.Range("T3").Value = 0.002382060785288784

$ws.Range("G4").Value = 12.96699333333333
$ws.Range("H4").Value = 38.90098
$ws.Range("I4").Value = 0.248939824633031
$ws.Range("J4").Value = 0.248939824633031
$ws.Range("M4").Value = 1.599526666666667
$ws.Range("N4").Value = 4.79858
$ws.Range("O4").Value = 0.008258927804608534
$ws.Range("P4").Value = 0.008258927804608534
$ws.Range("Q4").Value = 20.74105162315556
$ws.Range("R4").Value = 186.6694646084
$ws.Range("S4").Value = 0.002055976039336112
$ws.Range("T4").Value = 0.002055976039336112

$ws.Range("G5").Value = 12.96699333333333
$ws.Range("H5").Value = 38.90098
$ws.Range("I5").Value = 0.248939824633031
$ws.Range("J5").Value = 0.248939824633031
$ws.Range("M5").Value = 189.8529766666667
$ws.Range("N5").Value = 569.55893
$ws.Range("O5").Value = 0.9802787664976067
$ws.Range("P5").Value = 0.9802787664976066
$ws.Range("Q5").Value = 2461.822282750155
$ws.Range("R5").Value = 22156.4005447514
$ws.Range("S5").Value = 0.2440304242233981
$ws.Range("T5").Value = 0.2440304242233981

$ws.Range("I6").Value = 0.5739110916856193
$ws.Range("J6").Value = 0.5739110916856194
$ws.Range("M6").Value = 0.3667156666666667
$ws.Range("N6").Value = 1.100147
$ws.Range("O6").Value = 0.001893484040582144
$ws.Range("P6").Value = 0.001893484040582144
$ws.Range("Q6").Value = 10.96273688009389
$ws.Range("R6").Value = 98.664631920845
$ws.Range("S6").Value = 0.001086691492819795
$ws.Range("T6").Value = 0.001086691492819796

$ws.Range("I7").Value = 0.5739110916856193
$ws.Range("J7").Value = 0.5739110916856194
$ws.Range("N7").Value = 5.559651000000001
$ws.Range("O7").Value = 0.009568821657202681
$ws.Range("P7").Value = 0.009568821657202679
$ws.Range("Q7").Value = 55.400770131765
$ws.Range("R7").Value = 498.606931185885
$ws.Range("S7").Value = 0.005491652883430187
$ws.Range("T7").Value = 0.005491652883430187

$ws.Range("I8").Value = 0.5739110916856193
$ws.Range("J8").Value = 0.5739110916856194
$ws.Range("M8").Value = 1.599526666666667
$ws.Range("N8").Value = 4.79858
$ws.Range("O8").Value = 0.008258927804608534
$ws.Range("P8").Value = 0.008258927804608534
$ws.Range("Q8").Value = 47.81685532758889
$ws.Range("R8").Value = 430.3516979483
$ws.Range("S8").Value = 0.004739890272495598
$ws.Range("T8").Value = 0.004739890272495599

$ws.Range("I9").Value = 0.5739110916856193
$ws.Range("J9").Value = 0.5739110916856194
$ws.Range("M9").Value = 189.8529766666667
$ws.Range("N9").Value = 569.55893
$ws.Range("O9").Value = 0.9802787664976067
$ws.Range("P9").Value = 0.9802787664976066
$ws.Range("Q9").Value = 5675.536712182839
$ws.Range("R9").Value = 51079.83040964555
$ws.Range("S9").Value = 0.5625928570368737
$ws.Range("T9").Value = 0.5625928570368738

$ws.Range("G10").Value = 7.679779666666666
$ws.Range("H10").Value = 23.039339
$ws.Range("I10").Value = 0.1474361059880998
$ws.Range("J10").Value = 0.1474361059880998
$ws.Range("M10").Value = 0.3667156666666667
$ws.Range("N10").Value = 1.100147
$ws.Range("O10").Value = 0.001893484040582144
$ws.Range("P10").Value = 0.001893484040582144
$ws.Range("Q10").Value = 2.816295520314777
$ws.Range("R10").Value = 25.346659682833
$ws.Range("S10").Value = 0.0002791679136940444
$ws.Range("T10").Value = 0.0002791679136940444

$ws.Range("G11").Value = 7.679779666666666
$ws.Range("H11").Value = 23.039339
$ws.Range("I11").Value = 0.1474361059880998
$ws.Range("J11").Value = 0.1474361059880998
$ws.Range("N11").Value = 5.559651000000001
$ws.Range("O11").Value = 0.009568821657202681
$ws.Range("P11").Value = 0.009568821657202679
$ws.Range("Q11").Value = 14.232298234521
$ws.Range("R11").Value = 128.090684110689
$ws.Range("S11").Value = 0.001410789804032559
$ws.Range("T11").Value = 0.001410789804032559

$ws.Range("G12").Value = 7.679779666666666
$ws.Range("H12").Value = 23.039339
$ws.Range("I12").Value = 0.1474361059880998
$ws.Range("J12").Value = 0.1474361059880998
$ws.Range("M12").Value = 1.599526666666667
$ws.Range("N12").Value = 4.79858
$ws.Range("O12").Value = 0.008258927804608534
$ws.Range("P12").Value = 0.008258927804608534
$ws.Range("Q12").Value = 12.28401237095778
$ws.Range("R12").Value = 110.55611133862
$ws.Range("S12").Value = 0.001217664155148328
$ws.Range("T12").Value = 0.001217664155148328

$ws.Range("G13").Value = 7.679779666666666
$ws.Range("H13").Value = 23.039339
$ws.Range("I13").Value = 0.1474361059880998
$ws.Range("J13").Value = 0.1474361059880998
$ws.Range("M13").Value = 189.8529766666667
$ws.Range("N13").Value = 569.55893
$ws.Range("O13").Value = 0.9802787664976067
$ws.Range("P13").Value = 0.9802787664976066
$ws.Range("Q13").Value = 1458.029029860808
$ws.Range("R13").Value = 13122.26126874727
$ws.Range("S13").Value = 0.1445284841152249
$ws.Range("T13").Value = 0.1445284841152249

$ws.Range("G14").Value = 1.547715333333333
$ws.Range("H14").Value = 4.643146
$ws.Range("I14").Value = 0.02971297769324987
$ws.Range("J14").Value = 0.02971297769324987
$ws.Range("M14").Value = 0.3667156666666667
$ws.Range("N14").Value = 1.100147
$ws.Range("O14").Value = 0.001893484040582144
$ws.Range("P14").Value = 0.001893484040582144
$ws.Range("Q14").Value = 0.5675714602735555
$ws.Range("R14").Value = 5.108143142462
$ws.Range("S14").Value = 0.00005626104906034186
$ws.Range("T14").Value = 0.00005626104906034187

$ws.Range("G15").Value = 1.547715333333333
$ws.Range("H15").Value = 4.643146
$ws.Range("I15").Value = 0.02971297769324987
$ws.Range("J15").Value = 0.02971297769324987
$ws.Range("N15").Value = 5.559651000000001
$ws.Range("O15").Value = 0.009568821657202681
$ws.Range("P15").Value = 0.009568821657202679
$ws.Range("Q15").Value = 2.868252366894
$ws.Range("R15").Value = 25.814271302046
$ws.Range("S15").Value = 0.0002843181844511495
$ws.Range("T15").Value = 0.0002843181844511495

$ws.Range("G16").Value = 1.547715333333333
$ws.Range("H16").Value = 4.643146
$ws.Range("I16").Value = 0.02971297769324987
$ws.Range("J16").Value = 0.02971297769324987
$ws.Range("M16").Value = 1.599526666666667
$ws.Range("N16").Value = 4.79858
$ws.Range("O16").Value = 0.008258927804608534
$ws.Range("P16").Value = 0.008258927804608534
$ws.Range("Q16").Value = 2.475611948075556
$ws.Range("R16").Value = 22.28050753268
$ws.Range("S16").Value = 0.0002453973376284945
$ws.Range("T16").Value = 0.0002453973376284945

$ws.Range("G17").Value = 1.547715333333333
$ws.Range("H17").Value = 4.643146
$ws.Range("I17").Value = 0.02971297769324987
$ws.Range("J17").Value = 0.02971297769324987
$ws.Range("M17").Value = 189.8529766666667
$ws.Range("N17").Value = 569.55893
$ws.Range("O17").Value = 0.9802787664976067
$ws.Range("P17").Value = 0.9802787664976066
$ws.Range("Q17").Value = 293.8383630659756
$ws.Range("R17").Value = 2644.54526759378
$ws.Range("S17").Value = 0.02912700112210988
$ws.Range("T17").Value = 0.02912700112210988
